{"js": "// Replace the date line and the 25 two-digit-by-two-digit multiplication\n// equations in the table with the new values from the latest generated\n// output, while keeping all existing formatting (fonts, sizes, alignment).\nconst replacements = [\n  [\"2025-03-06 Thursday\", \"2025-03-07 Friday\"],\n  [\"17\u00d765=1105\", \"83\u00d761=5063\"],\n  [\"12\u00d799=1188\", \"61\u00d750=3050\"],\n  [\"66\u00d781=5346\", \"29\u00d752=1508\"],\n  [\"71\u00d711=781\", \"97\u00d713=1261\"],\n  [\"43\u00d715=645\", \"61\u00d758=3538\"],\n  [\"57\u00d782=4674\", \"17\u00d750=850\"],\n  [\"35\u00d758=2030\", \"59\u00d775=4425\"],\n  [\"39\u00d772=2808\", \"26\u00d790=2340\"],\n  [\"38\u00d765=2470\", \"42\u00d738=1596\"],\n  [\"51\u00d786=4386\", \"53\u00d714=742\"],\n  [\"44\u00d720=880\", \"52\u00d766=3432\"],\n  [\"33\u00d729=957\", \"25\u00d748=1200\"],\n  [\"47\u00d717=799\", \"26\u00d734=884\"],\n  [\"28\u00d754=1512\", \"98\u00d792=9016\"],\n  [\"86\u00d758=4988\", \"45\u00d785=3825\"],\n  [\"25\u00d760=1500\", \"51\u00d744=2244\"],\n  [\"11\u00d797=1067\", \"75\u00d761=4575\"],\n  [\"84\u00d777=6468\", \"34\u00d711=374\"],\n  [\"82\u00d780=6560\", \"62\u00d736=2232\"],\n  [\"29\u00d770=2030\", \"85\u00d790=7650\"],\n  [\"21\u00d778=1638\", \"60\u00d726=1560\"],\n  [\"49\u00d772=3528\", \"51\u00d790=4590\"],\n  [\"70\u00d753=3710\", \"12\u00d753=636\"],\n  [\"50\u00d768=3400\", \"28\u00d798=2744\"],\n  [\"55\u00d786=4730\", \"16\u00d792=1472\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 two-digit-by-two-digit multiplication\n# equations in the table with the new values from the latest generated\n# output, while keeping all existing formatting (fonts, sizes, alignment).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-03-06 Thursday\", \"2025-03-07 Friday\"),\n    @(\"17\u00d765=1105\", \"83\u00d761=5063\"),\n    @(\"12\u00d799=1188\", \"61\u00d750=3050\"),\n    @(\"66\u00d781=5346\", \"29\u00d752=1508\"),\n    @(\"71\u00d711=781\", \"97\u00d713=1261\"),\n    @(\"43\u00d715=645\", \"61\u00d758=3538\"),\n    @(\"57\u00d782=4674\", \"17\u00d750=850\"),\n    @(\"35\u00d758=2030\", \"59\u00d775=4425\"),\n    @(\"39\u00d772=2808\", \"26\u00d790=2340\"),\n    @(\"38\u00d765=2470\", \"42\u00d738=1596\"),\n    @(\"51\u00d786=4386\", \"53\u00d714=742\"),\n    @(\"44\u00d720=880\", \"52\u00d766=3432\"),\n    @(\"33\u00d729=957\", \"25\u00d748=1200\"),\n    @(\"47\u00d717=799\", \"26\u00d734=884\"),\n    @(\"28\u00d754=1512\", \"98\u00d792=9016\"),\n    @(\"86\u00d758=4988\", \"45\u00d785=3825\"),\n    @(\"25\u00d760=1500\", \"51\u00d744=2244\"),\n    @(\"11\u00d797=1067\", \"75\u00d761=4575\"),\n    @(\"84\u00d777=6468\", \"34\u00d711=374\"),\n    @(\"82\u00d780=6560\", \"62\u00d736=2232\"),\n    @(\"29\u00d770=2030\", \"85\u00d790=7650\"),\n    @(\"21\u00d778=1638\", \"60\u00d726=1560\"),\n    @(\"49\u00d772=3528\", \"51\u00d790=4590\"),\n    @(\"70\u00d753=3710\", \"12\u00d753=636\"),\n    @(\"50\u00d768=3400\", \"28\u00d798=2744\"),\n    @(\"55\u00d786=4730\", \"16\u00d792=1472\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
